# Implement loading/saving of error values.
#
# 1. "Cell Values" sheet: add a row of #DIV/0! error values (row 8),
#    mirroring the existing layout (B:E hold the raw value, F/G hold the
#    string representation used by GetString()/GetFormattedString()).
# 2. Add a new "Errors" worksheet (after "Test Whitespace") demonstrating
#    every Excel error literal plus a formula that produces the same error.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- "Cell Values" sheet: new row 8 with #DIV/0! values ---
$ws1.Range("B8").Value = "#DIV/0!"
$ws1.Range("C8").Value = "#DIV/0!"
$ws1.Range("D8").Value = "#DIV/0!"
$ws1.Range("E8").Value = "#DIV/0!"
# Leading apostrophe -> literal text "#DIV/0!" (quote-prefixed), not an error.
$ws1.Range("F8").Value = "'#DIV/0!"
$ws1.Range("G8").Value = "'#DIV/0!"

# --- New "Errors" worksheet, placed after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$errorsSheet = $wb.Worksheets.Add($null, $lastSheet)
$errorsSheet.Name = "Errors"

$errorsSheet.Range("B2").Value = "Error value"
$errorsSheet.Range("C2").Value = "Formula error"

$errorsSheet.Range("B3").Value = "#REF!"
$errorsSheet.Range("C3").Formula = "=#REF!+1"

$errorsSheet.Range("B4").Value = "#VALUE!"
$errorsSheet.Range("C4").Formula = "=`"TRUE`"*1"

$errorsSheet.Range("B5").Value = "#DIV/0!"
$errorsSheet.Range("C5").Formula = "=1/0"

$errorsSheet.Range("B6").Value = "#NAME?"
$errorsSheet.Range("C6").Formula = "=NONEXISTENT.FUNCTION()"

$errorsSheet.Range("B7").Value = "#N/A"
$errorsSheet.Range("C7").Formula = "=NA()"

$errorsSheet.Range("B8").Value = "#NULL!"
$errorsSheet.Range("C8").Formula = "=#NULL!+1"

$errorsSheet.Range("B9").Value = "#NUM!"
$errorsSheet.Range("C9").Formula = "=#NUM!+1"
